$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.718438863754272
$ws.Range("B1").Value = 2.300286054611206
$ws.Range("C1").Value = 3.355501413345337
$ws.Range("D1").Value = 4.116728782653809
$ws.Range("E1").Value = 0.6514696478843689
